$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'29.902.28"
$ws.Range("D2").Style = "Normal"
$ws.Range("D3").Value = "'1.892.84"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  -0.05%  "
$ws.Range("E3").Style = "Normal"
$ws.Range("E4").Value = "'  -0.06%  "
$ws.Range("E4").Style = "Normal"
$ws.Range("E5").Value = "'  -1.97%  "
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'244.05"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'  +0.20%  "
$ws.Range("E6").Style = "Normal"
$ws.Range("E7").Value = "'  -0.01%  "
$ws.Range("E7").Style = "Normal"
$ws.Range("D8").Value = "'0.3126"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'  -0.67%  "
$ws.Range("E8").Style = "Normal"
$ws.Range("D9").Value = "'25.69"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'  +1.55%  "
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Value = "'0.07232"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'  +0.39%  "
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Value = "'0.08693"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'  +7.52%  "
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = "'2.043.68"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'  +7.86%  "
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = "'0.7711"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'  +0.78%  "
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = "'5.401"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'  -1.79%  "
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'94.23"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'  +1.92%  "
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = "'6.200"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'  +0.88%  "
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = "'30.122.10"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'  +0.91%  "
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = "'13.92"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'  +0.01%  "
$ws.Range("E18").Style = "Normal"
$ws.Range("B19").Value = "'WrappedliquidstakedEther2.0"
$ws.Range("B19").Style = "Normal"
$ws.Range("C19").Value = "'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("C19").Style = "Normal"
$ws.Range("D19").Value = "'2.307.72"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'  +8.07%  "
$ws.Range("E19").Style = "Normal"
$ws.Range("B20").Value = "'BitcoinCash"
$ws.Range("B20").Style = "Normal"
$ws.Range("C20").Value = "'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("C20").Style = "Normal"
$ws.Range("D20").Value = "'245.25"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'  +0.61%  "
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = "'0.000007858"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'  +1.01%  "
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = "'8.164"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'  +0.66%  "
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = "'1.001"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'  +0.03%  "
$ws.Range("E23").Style = "Normal"
$ws.Range("E24").Value = "'  -0.12%  "
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = "'0.1592"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'  -3.23%  "
$ws.Range("E25").Style = "Normal"
$ws.Range("D26").Value = "'9.520"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'  +1.21%  "
$ws.Range("E26").Style = "Normal"
$ws.Range("D27").Value = "'162.38"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "'  -0.58%  "
$ws.Range("E27").Style = "Normal"
$ws.Range("E28").Value = "'  +0.41%  "
$ws.Range("E28").Style = "Normal"
$ws.Range("D29").Value = "'2.041"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "'  -0.52%  "
$ws.Range("E29").Style = "Normal"
$ws.Range("D30").Value = "'1.434"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "'  +1.80%  "
$ws.Range("E30").Style = "Normal"
$ws.Range("D31").Value = "'1.544"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "'  -0.31%  "
$ws.Range("E31").Style = "Normal"
$ws.Range("D32").Value = "'4.516"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "'  +0.17%  "
$ws.Range("E32").Style = "Normal"
$ws.Range("D33").Value = "'4.115"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "'  +0.05%  "
$ws.Range("E33").Style = "Normal"
$ws.Range("D34").Value = "'0.05467"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "'  -1.67%  "
$ws.Range("E34").Style = "Normal"
$ws.Range("E35").Value = "'  -1.77%  "
$ws.Range("E35").Style = "Normal"
$ws.Range("D36").Value = "'0.7532"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "'  +1.43%  "
$ws.Range("E36").Style = "Normal"
$ws.Range("D37").Value = "'1.007"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "'  +0.99%  "
$ws.Range("E37").Style = "Normal"
$ws.Range("D38").Value = "'2.701"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "'  +3.20%  "
$ws.Range("E38").Style = "Normal"
$ws.Range("E39").Value = "'  +2.61%  "
$ws.Range("E39").Style = "Normal"
$ws.Range("E40").Value = "'  +0.32%  "
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = "'0.4509"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'  +2.07%  "
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Value = "'73.59"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'  -0.51%  "
$ws.Range("E42").Style = "Normal"
$ws.Range("B43").Value = "'FraxShare"
$ws.Range("B43").Style = "Normal"
$ws.Range("C43").Value = "'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("C43").Style = "Normal"
$ws.Range("D43").Value = "'6.082"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'  +3.73%  "
$ws.Range("E43").Style = "Normal"
$ws.Range("B44").Value = "'Maker"
$ws.Range("B44").Style = "Normal"
$ws.Range("C44").Value = "'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("C44").Style = "Normal"
$ws.Range("D44").Value = "'1.096.35"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'  -3.93%  "
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = "'0.8543"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'  +0.54%  "
$ws.Range("E45").Style = "Normal"
$ws.Range("D46").Value = "'2.217.89"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'  +8.33%  "
$ws.Range("E46").Style = "Normal"
$ws.Range("E47").Value = "'  -0.05%  "
$ws.Range("E47").Style = "Normal"
$ws.Range("D48").Value = "'103.05"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'  -0.79%  "
$ws.Range("E48").Style = "Normal"
$ws.Range("D50").Value = "'7.618"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'  +2.12%  "
$ws.Range("E50").Style = "Normal"
$ws.Range("D51").Value = "'9.856"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'  -1.24%  "
$ws.Range("E51").Style = "Normal"
